# close #187: Remove unnecessary column name in values and proportionality
# The "nome" column (column B) is no longer needed; delete it so that all
# the remaining columns shift one position to the left.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire column B, shifting C:M left to B:L
$ws.Columns.Item(2).Delete()

# Move selection to B1 (matches the resulting saved selection in the sheet)
[void]$ws.Range("B1").Select()
